$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 206, shifting existing rows 206-254 down to 207-255.
$ws.Rows.Item(206).EntireRow.Insert()

# Populate the newly inserted row 206 with the new data record.
$ws.Range("A206").Value = 10
$ws.Range("B206").Value = 'Vega Modelo de Temuco'
$ws.Range("C206").Value = 'La Araucanía'
$ws.Range("D206").Value = 44642
$ws.Range("E206").Value = 9
$ws.Range("F206").Value = 100112001
$ws.Range("G206").Value = 'Berenjena'
$ws.Range("H206").Value = 'Sin especificar'
$ws.Range("I206").Value = 'Primera'
$ws.Range("J206").Value = 45
$ws.Range("K206").Value = 12000
$ws.Range("L206").Value = 13000
$ws.Range("M206").Value = 12556
$ws.Range("N206").Value = '$/caja 60 unidades'
$ws.Range("O206").Value = 'Región del Maule'
$ws.Range("P206").Value = 209
$ws.Range("Q206").Value = 60
$ws.Range("R206").Value = 'Hortaliza'
